# Supplemental Table 1 — update with review process.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove all the leftover hidden "_xlchart.*" defined names (chart data
#    links from a deleted chart) that cluttered the workbook.
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# 2. Fix the placeholder table caption: "Table Y" -> "Table 1".
$ws.Range("A1").Value = "Supplemental Table 1. Average pH associated with overlying water of each wetland region. Avg pH = average pH, STDev = standard deviation, n = number of replicate readings taken per region."
